$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12: add response_collected (column E) value
$ws.Range("E12").Value = "master_worker_response_tracke_Sep-22-2023.csv"

# Row 13: new row of data - order matches shared string insertion order
$ws.Range("E13").Value = "master_worker_response_tracke_resub_Sep-22-2023.csv"
$ws.Range("A13").Value = "paiewise_resub"
$ws.Range("D13").Value = "all_submitted_tracker_nina_resubSep-22-2023.csv"
$ws.Range("F13").Value = "master_all_responses_Sep-22-2023_to_resub_Sep-22-2023_Nina.csv"
# Copy B12 ("Sep-22-2023", stored as text) into B13 to avoid Excel's
# auto date-conversion of a plain string assignment
$ws.Range("B12").Copy($ws.Range("B13"))
$ws.Range("C13").Value = "NV"

# Update selection to match
$ws.Range("E15").Select()
